# Auto-generated edit script: updates TPM-derived NATMI metric columns
# (G,H,I,J,M,N,O,P,Q,R,S,T) for rows 2-26 per the commit "update scripts wuth new tpm"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 32.30926933333333
$ws.Range("H2").Value = 96.927808
$ws.Range("I2").Value = 0.4126751607889972
$ws.Range("J2").Value = 0.4140835826598944
$ws.Range("M2").Value = 293.7577056666667
$ws.Range("N2").Value = 881.273117
$ws.Range("O2").Value = 0.9369756110667984
$ws.Range("P2").Value = 0.940660486426629
$ws.Range("Q2").Value = 9491.096831126391
$ws.Range("R2").Value = 85419.87148013753
$ws.Range("S2").Value = 0.3866665609523599
$ws.Range("T2").Value = 0.3895120642861375
$ws.Range("G3").Value = 32.30926933333333
$ws.Range("H3").Value = 96.927808
$ws.Range("I3").Value = 0.4126751607889972
$ws.Range("J3").Value = 0.4140835826598944
$ws.Range("O3").Value = 0.01525159481997056
$ws.Range("P3").Value = 0.01531157527761154
$ws.Range("Q3").Value = 154.4910684501547
$ws.Range("R3").Value = 1390.419616051392
$ws.Range("S3").Value = 0.006293954344619986
$ws.Range("T3").Value = 0.006340271947120055
$ws.Range("G4").Value = 32.30926933333333
$ws.Range("H4").Value = 96.927808
$ws.Range("I4").Value = 0.4126751607889972
$ws.Range("J4").Value = 0.4140835826598944
$ws.Range("M4").Value = 7.906212666666666
$ws.Range("N4").Value = 23.718638
$ws.Range("O4").Value = 0.02521781829607561
$ws.Range("P4").Value = 0.02531699325449539
$ws.Range("Q4").Value = 255.4439544539449
$ws.Range("R4").Value = 2298.995590085504
$ws.Range("S4").Value = 0.01040676722008072
$ws.Range("T4").Value = 0.01048335126899783
$ws.Range("G5").Value = 32.30926933333333
$ws.Range("H5").Value = 96.927808
$ws.Range("I5").Value = 0.4126751607889972
$ws.Range("J5").Value = 0.4140835826598944
$ws.Range("M5").Value = 3.6844455
$ws.Range("N5").Value = 7.368891000000001
$ws.Range("O5").Value = 0.01175198303639443
$ws.Range("P5").Value = 0.007865466969060864
$ws.Range("Q5").Value = 119.041742003488
$ws.Range("R5").Value = 714.250452020928
$ws.Range("S5").Value = 0.004849751489133641
$ws.Range("T5").Value = 0.003256960741841783
$ws.Range("G6").Value = 32.30926933333333
$ws.Range("H6").Value = 96.927808
$ws.Range("I6").Value = 0.4126751607889972
$ws.Range("J6").Value = 0.4140835826598944
$ws.Range("M6").Value = 3.386921
$ws.Range("N6").Value = 10.160763
$ws.Range("O6").Value = 0.01080299278076119
$ws.Range("P6").Value = 0.01084547807220323
$ws.Range("Q6").Value = 109.4289427997226
$ws.Range("R6").Value = 984.860485197504
$ws.Range("S6").Value = 0.004458126782803002
$ws.Range("T6").Value = 0.00449093441579724
$ws.Range("I7").Value = 0.03337290046497914
$ws.Range("J7").Value = 0.03348679906459509
$ws.Range("M7").Value = 293.7577056666667
$ws.Range("N7").Value = 881.273117
$ws.Range("O7").Value = 0.9369756110667984
$ws.Range("P7").Value = 0.940660486426629
$ws.Range("Q7").Value = 767.541785754858
$ws.Range("R7").Value = 6907.876071793722
$ws.Range("S7").Value = 0.03126959380624527
$ws.Range("T7").Value = 0.03149970869697281
$ws.Range("I8").Value = 0.03337290046497914
$ws.Range("J8").Value = 0.03348679906459509
$ws.Range("O8").Value = 0.01525159481997056
$ws.Range("P8").Value = 0.01531157527761154
$ws.Range("S8").Value = 0.0005089899558590689
$ws.Range("T8").Value = 0.0005127356446837996
$ws.Range("I9").Value = 0.03337290046497914
$ws.Range("J9").Value = 0.03348679906459509
$ws.Range("M9").Value = 7.906212666666666
$ws.Range("N9").Value = 23.718638
$ws.Range("O9").Value = 0.02521781829607561
$ws.Range("P9").Value = 0.02531699325449539
$ws.Range("Q9").Value = 20.65766606856911
$ws.Range("R9").Value = 185.918994617122
$ws.Range("S9").Value = 0.0008415917399388613
$ws.Range("T9").Value = 0.0008477850660329964
$ws.Range("I10").Value = 0.03337290046497914
$ws.Range("J10").Value = 0.03348679906459509
$ws.Range("M10").Value = 3.6844455
$ws.Range("N10").Value = 7.368891000000001
$ws.Range("O10").Value = 0.01175198303639443
$ws.Range("P10").Value = 0.007865466969060864
$ws.Range("Q10").Value = 9.6268653520715
$ws.Range("R10").Value = 57.761192112429
$ws.Range("S10").Value = 0.0003921977601397148
$ws.Range("T10").Value = 0.0002633893119421509
$ws.Range("I11").Value = 0.03337290046497914
$ws.Range("J11").Value = 0.03348679906459509
$ws.Range("M11").Value = 3.386921
$ws.Range("N11").Value = 10.160763
$ws.Range("O11").Value = 0.01080299278076119
$ws.Range("P11").Value = 0.01084547807220323
$ws.Range("Q11").Value = 8.849481536666332
$ws.Range("R11").Value = 79.64533382999699
$ws.Range("S11").Value = 0.0003605272027962316
$ws.Range("T11").Value = 0.0003631803449633418
$ws.Range("G12").Value = 12.89411533333333
$ws.Range("H12").Value = 38.682346
$ws.Range("I12").Value = 0.1646920907903501
$ws.Range("J12").Value = 0.1652541695502867
$ws.Range("M12").Value = 293.7577056666667
$ws.Range("N12").Value = 881.273117
$ws.Range("O12").Value = 0.9369756110667984
$ws.Range("P12").Value = 0.940660486426629
$ws.Range("Q12").Value = 3787.745736921386
$ws.Range("R12").Value = 34089.71163229248
$ws.Range("S12").Value = 0.1543124724061569
$ws.Range("T12").Value = 0.1554480675132013
$ws.Range("G13").Value = 12.89411533333333
$ws.Range("H13").Value = 38.682346
$ws.Range("I13").Value = 0.1646920907903501
$ws.Range("J13").Value = 0.1652541695502867
$ws.Range("O13").Value = 0.01525159481997056
$ws.Range("P13").Value = 0.01531157527761154
$ws.Range("Q13").Value = 61.65492738367266
$ws.Range("R13").Value = 554.8943464530539
$ws.Range("S13").Value = 0.002511817038788224
$ws.Range("T13").Value = 0.002530301657008396
$ws.Range("G14").Value = 12.89411533333333
$ws.Range("H14").Value = 38.682346
$ws.Range("I14").Value = 0.1646920907903501
$ws.Range("J14").Value = 0.1652541695502867
$ws.Range("M14").Value = 7.906212666666666
$ws.Range("N14").Value = 23.718638
$ws.Range("O14").Value = 0.02521781829607561
$ws.Range("P14").Value = 0.02531699325449539
$ws.Range("Q14").Value = 101.9436179738609
$ws.Range("R14").Value = 917.4925617647478
$ws.Range("S14").Value = 0.004153175220351836
$ws.Range("T14").Value = 0.004183738695781846
$ws.Range("G15").Value = 12.89411533333333
$ws.Range("H15").Value = 38.682346
$ws.Range("I15").Value = 0.1646920907903501
$ws.Range("J15").Value = 0.1652541695502867
$ws.Range("M15").Value = 3.6844455
$ws.Range("N15").Value = 7.368891000000001
$ws.Range("O15").Value = 0.01175198303639443
$ws.Range("P15").Value = 0.007865466969060864
$ws.Range("Q15").Value = 47.507665216381
$ws.Range("R15").Value = 285.045991298286
$ws.Range("S15").Value = 0.001935458657196526
$ws.Range("T15").Value = 0.001299801212097363
$ws.Range("G16").Value = 12.89411533333333
$ws.Range("H16").Value = 38.682346
$ws.Range("I16").Value = 0.1646920907903501
$ws.Range("J16").Value = 0.1652541695502867
$ws.Range("M16").Value = 3.386921
$ws.Range("N16").Value = 10.160763
$ws.Range("O16").Value = 0.01080299278076119
$ws.Range("P16").Value = 0.01084547807220323
$ws.Range("Q16").Value = 43.67134999888866
$ws.Range("R16").Value = 393.0421499899979
$ws.Range("S16").Value = 0.001779167467856619
$ws.Range("T16").Value = 0.001792260472197789
$ws.Range("G17").Value = 0.798886
$ws.Range("H17").Value = 1.597772
$ws.Range("I17").Value = 0.01020389551681842
$ws.Range("J17").Value = 0.006825813640948785
$ws.Range("M17").Value = 293.7577056666667
$ws.Range("N17").Value = 881.273117
$ws.Range("O17").Value = 0.9369756110667984
$ws.Range("P17").Value = 0.940660486426629
$ws.Range("Q17").Value = 234.6789184492206
$ws.Range("R17").Value = 1408.073510695324
$ws.Range("S17").Value = 0.009560801237132703
$ws.Range("T17").Value = 0.006420773179752403
$ws.Range("G18").Value = 0.798886
$ws.Range("H18").Value = 1.597772
$ws.Range("I18").Value = 0.01020389551681842
$ws.Range("J18").Value = 0.006825813640948785
$ws.Range("O18").Value = 0.01525159481997056
$ws.Range("P18").Value = 0.01531157527761154
$ws.Range("Q18").Value = 3.819979660838
$ws.Range("R18").Value = 22.919877965028
$ws.Range("S18").Value = 0.0001556256800078286
$ws.Range("T18").Value = 0.000104513959394335
$ws.Range("G19").Value = 0.798886
$ws.Range("H19").Value = 1.597772
$ws.Range("I19").Value = 0.01020389551681842
$ws.Range("J19").Value = 0.006825813640948785
$ws.Range("M19").Value = 7.906212666666666
$ws.Range("N19").Value = 23.718638
$ws.Range("O19").Value = 0.02521781829607561
$ws.Range("P19").Value = 0.02531699325449539
$ws.Range("Q19").Value = 6.316162612422666
$ws.Range("R19").Value = 37.89697567453599
$ws.Range("S19").Value = 0.0002573199830552674
$ws.Range("T19").Value = 0.000172809077904343
$ws.Range("G20").Value = 0.798886
$ws.Range("H20").Value = 1.597772
$ws.Range("I20").Value = 0.01020389551681842
$ws.Range("J20").Value = 0.006825813640948785
$ws.Range("M20").Value = 3.6844455
$ws.Range("N20").Value = 7.368891000000001
$ws.Range("O20").Value = 0.01175198303639443
$ws.Range("P20").Value = 0.007865466969060864
$ws.Range("Q20").Value = 2.943451927713
$ws.Range("R20").Value = 11.773807710852
$ws.Range("S20").Value = 0.0001199160070187913
$ws.Range("T20").Value = 0.00005368821172984774
$ws.Range("G21").Value = 0.798886
$ws.Range("H21").Value = 1.597772
$ws.Range("I21").Value = 0.01020389551681842
$ws.Range("J21").Value = 0.006825813640948785
$ws.Range("M21").Value = 3.386921
$ws.Range("N21").Value = 10.160763
$ws.Range("O21").Value = 0.01080299278076119
$ws.Range("P21").Value = 0.01084547807220323
$ws.Range("Q21").Value = 2.705763770006
$ws.Range("R21").Value = 16.234582620036
$ws.Range("S21").Value = 0.0001102326096038309
$ws.Range("T21").Value = 0.00007402921216785576
$ws.Range("G22").Value = 29.67714566666667
$ws.Range("H22").Value = 89.03143700000001
$ws.Range("I22").Value = 0.3790559524388551
$ws.Range("J22").Value = 0.3803496350842752
$ws.Range("M22").Value = 293.7577056666667
$ws.Range("N22").Value = 881.273117
$ws.Range("O22").Value = 0.9369756110667984
$ws.Range("P22").Value = 0.940660486426629
$ws.Range("Q22").Value = 8717.89022177546
$ws.Range("R22").Value = 78461.01199597913
$ws.Range("S22").Value = 0.3551661826649035
$ws.Range("T22").Value = 0.3577798727505651
$ws.Range("G23").Value = 29.67714566666667
$ws.Range("H23").Value = 89.03143700000001
$ws.Range("I23").Value = 0.3790559524388551
$ws.Range("J23").Value = 0.3803496350842752
$ws.Range("O23").Value = 0.01525159481997056
$ws.Range("P23").Value = 0.01531157527761154
$ws.Range("Q23").Value = 141.9052190655404
$ws.Range("R23").Value = 1277.146971589863
$ws.Range("S23").Value = 0.005781207800695448
$ws.Range("T23").Value = 0.00582375206940496
$ws.Range("G24").Value = 29.67714566666667
$ws.Range("H24").Value = 89.03143700000001
$ws.Range("I24").Value = 0.3790559524388551
$ws.Range("J24").Value = 0.3803496350842752
$ws.Range("M24").Value = 7.906212666666666
$ws.Range("N24").Value = 23.718638
$ws.Range("O24").Value = 0.02521781829607561
$ws.Range("P24").Value = 0.02531699325449539
$ws.Range("Q24").Value = 234.6338249803118
$ws.Range("R24").Value = 2111.704424822806
$ws.Range("S24").Value = 0.009558964132648926
$ws.Range("T24").Value = 0.009629309145778378
$ws.Range("G25").Value = 29.67714566666667
$ws.Range("H25").Value = 89.03143700000001
$ws.Range("I25").Value = 0.3790559524388551
$ws.Range("J25").Value = 0.3803496350842752
$ws.Range("M25").Value = 3.6844455
$ws.Range("N25").Value = 7.368891000000001
$ws.Range("O25").Value = 0.01175198303639443
$ws.Range("P25").Value = 0.007865466969060864
$ws.Range("Q25").Value = 109.3438258043945
$ws.Range("R25").Value = 656.0629548263671
$ws.Range("S25").Value = 0.004454659122905761
$ws.Range("T25").Value = 0.00299162749144972
$ws.Range("G26").Value = 29.67714566666667
$ws.Range("H26").Value = 89.03143700000001
$ws.Range("I26").Value = 0.3790559524388551
$ws.Range("J26").Value = 0.3803496350842752
$ws.Range("M26").Value = 3.386921
$ws.Range("N26").Value = 10.160763
$ws.Range("O26").Value = 0.01080299278076119
$ws.Range("P26").Value = 0.01084547807220323
$ws.Range("Q26").Value = 100.5141478784923
$ws.Range("R26").Value = 904.6273309064311
$ws.Range("S26").Value = 0.004125073627077008
$ws.Range("T26").Value = 0.004125073627077008
